$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Fill in row 2 of Sheet3 (A2:F2) - mirrors row 1's table/column names, with a new
# combined "columns" value going into the (already wrap-text styled) C2/F2 cells
$ws3.Range("A2").Value = "org"
$ws3.Range("B2").Value = "employee"
$ws3.Range("C2").Value = "empno,ename,sal,job,doj,deptno,manager_id,bonus"
$ws3.Range("D2").Value = "org"
$ws3.Range("E2").Value = "employee"
$ws3.Range("F2").Value = "empno,ename,sal,job,doj,deptno,manager_id,bonus"

# Ensure the combined-columns cells wrap (C2/F2 already carried this style, but make
# it explicit/robust in case the pre-existing style isn't preserved on write)
$ws3.Range("C2").WrapText = $true
$ws3.Range("F2").WrapText = $true

# Row grew taller once it held the wrapped, comma-separated column list
$ws3.Rows.Item(2).RowHeight = 34

# Update selections / zoom to match author's final view
$ws1.Select()
$ws1.Range("B17").Select()

$ws3.Select()
$ws3.Range("E7").Select()

$win = $excel.ActiveWindow
$win.Zoom = 171

# Reposition/resize the app window to match the author's final screen layout
$win.Left = 12320
$win.Top = 3540
$win.Width = 28480
$win.Height = 14060

# Make Sheet3 the active sheet/tab (also leaves it selected/frontmost)
$ws3.Activate()
